# Deploy the implementation guide.
# Updates the CodeSystem metadata (Date/Count) and appends the newly
# published NCIT sample-type concepts to the Concepts sheet.

$wb = $excel.ActiveWorkbook

$meta = $wb.Worksheets.Item("Metadata")
$concepts = $wb.Worksheets.Item("Concepts")

# --- Metadata sheet: bump the publication Date and the concept Count. ---
$meta.Cells.Item(8, 2).Value = "2023-12-04T14:59:45+00:00"

# "Count" must stay a text value ("15"), not be auto-coerced to a number,
# while keeping the existing cell style. Build it as a text formula in a
# scratch cell, then paste-special just the resulting value over the
# target cell (this keeps the destination's original style untouched).
$scratch = $meta.Cells.Item(30, 10)
$scratch.Formula = "=""15"""
$scratch.Copy()
$countCell = $meta.Cells.Item(21, 2)
$countCell.PasteSpecial(-4163)  # xlPasteValues
$scratch.ClearContents()

# --- Concepts sheet: append the 9 newly introduced NCIT concepts. ---
# Copy row 7's formatting/content down into the new rows first so the new
# rows inherit the same styles (including the blank, styled Definition
# cell in column D) as the existing concept rows.
for ($r = 8; $r -le 16; $r++) {
    $srcRow = $concepts.Range("A7:D7")
    $destRow = $concepts.Range("A" + $r + ":D" + $r)
    $srcRow.Copy($destRow)
}

$newConcepts = @(
    @("NCIT:C13300", "Umbilical Cord Blood"),
    @("NCIT:C34320", "Umbilical Cord"),
    @("NCIT:C156445", "Derived Cell Line"),
    @("NCIT:156440", "Metastatic Tumor Sample"),
    @("NCIT:156441", "Sample Derived from New Primary"),
    @("NCIT:164032", "Tumor-Adjacent Normal Specimen"),
    @("NCIT:C18009", "Tumor Tissue"),
    @("NCIT:C162623", "Normal Tissue Segment"),
    @("NCIT:C156443", "Cell Line-Derived Xenograft")
)

$rowIndex = 8
foreach ($row in $newConcepts) {
    $concepts.Cells.Item($rowIndex, 2).Value = $row[0]
    $concepts.Cells.Item($rowIndex, 3).Value = $row[1]
    $rowIndex++
}
